# points_de_ramasse.xlsx - "some more work on new data"
#
# LECLERC ROUFFIAC no longer collects on Tuesdays ("Mardi") - update the
# "Jours de Ramasse" (pickup days) cell for that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = "Lundi, Mercredi, Vendredi"

# Row 1 (header) gets a touch more breathing room.
$ws.Rows(1).RowHeight = 20.25

# Normalize the font color used for the Nom/Adresse columns (A2:B8) so it
# matches the rest of the data rows (plain black) instead of the
# theme-derived color - this was an inconsistency left over from earlier
# edits.
$ws.Range("A2:B8").Font.Color = 0
